# Generate Report for Handoff
#
# A localization handoff report was regenerated: the four "Ready for
# handoff" rows (a73a235f…, ae17466c…, bbfb53fd…, c6626d81…) in both the
# "zh-cn" and "de-de" sheets now show Priority = "ht" (was "low"), and
# their "Latest Handoff Datetime" stamp was bumped by the new handoff run.

$wb = $excel.ActiveWorkbook

$rows = 4..7

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($r in $rows) {
        # Priority column (E): "low" -> "ht"
        $ws.Range("E$r").Value = "ht"
    }
}

# Latest Handoff Datetime column (H) for the same rows, per sheet.
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("H$r").Value = "2016-08-18 04:30:34"
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("H$r").Value = "2016-08-18 04:30:39"
}

# The "Overview" tab mirrors the de-de handoff timestamp in its own
# "Latest HO Xliff Generate Date" column (G) for these same four files, so
# it picks up the same refreshed stamp as part of this handoff run.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-18 04:30:39"
}
